# Falta arreglar cuando se da una fecha
#
# Two new batches of resume rows (6 entries each, for the 19:17 and 19:21
# runs on 06/11/2025) need to be appended to the data table, above the
# trailing "${table:resumeData.*}" placeholder row so the template row
# keeps acting as the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 blank rows right where the placeholder row currently sits.
# Excel shifts the placeholder row (and inherits formatting for the new
# rows) down by 12, so it ends up as the new last row.
$ws.Range("A56:A67").EntireRow.Insert()

# Pre-format the tradeDate column as text so date-looking strings like
# "06/11/2025" are stored verbatim instead of being auto-converted to a
# date serial number (the bug referenced in the commit message).
$ws.Range("B56:B67").NumberFormat = "@"

$newRows = @(
    @("2025-06-12 19:17:34", "06/11/2025", "CMEGroup Chicago-CU",   13),
    @("2025-06-12 19:17:35", "06/11/2025", "CMEGroup New York-NYH",  4),
    @("2025-06-12 19:17:36", "06/11/2025", "CMEGroup T2",           13),
    @("2025-06-12 19:17:36", "06/11/2025", "CMEGroup Corn",         15),
    @("2025-06-12 19:17:37", "06/11/2025", "CMEGroup RBob",         43),
    @("2025-06-12 19:17:37", "06/11/2025", "CMEGroup Sugar 11",      8),
    @("2025-06-12 19:21:26", "06/11/2025", "CMEGroup Chicago-CU",   13),
    @("2025-06-12 19:21:27", "06/11/2025", "CMEGroup New York-NYH",  4),
    @("2025-06-12 19:21:27", "06/11/2025", "CMEGroup T2",           13),
    @("2025-06-12 19:21:28", "06/11/2025", "CMEGroup Corn",         15),
    @("2025-06-12 19:21:28", "06/11/2025", "CMEGroup RBob",         43),
    @("2025-06-12 19:21:29", "06/11/2025", "CMEGroup Sugar 11",      8)
)

$row = 56
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

# Re-normalize the tradeDate cells' formatting to match the rest of the
# column (plain/general style) now that the text is safely stored.
$ws.Range("B55").Copy()
$ws.Range("B56:B67").PasteSpecial(-4122)
$excel.CutCopyMode = $false
